$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: empty "thick-bottom-border" cell, style matches the rest of row 3 (L3) ---
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

# --- Row 4: new "2022" column header, style matches the rest of row 4 (L4) ---
$ws.Range("M4").Value = 2022
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)

# --- Row 6: new data point, style matches the rest of row 6 (L6) ---
$ws.Range("M6").Value = 18
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)

# --- Row 7: new data point, style matches the rest of row 7 (L7) ---
$ws.Range("M7").Value = 6.2
$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)

# --- Row 8: new "-" placeholder, style matches the rest of row 8 (L8) ---
$ws.Range("M8").Value = "-"
$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)

# Clear clipboard marquee / move selection to match the saved workbook state
$excel.CutCopyMode = $false
$ws.Range("N4").Select() | Out-Null
